# Correcting names for outputs and obs data for soil mineral n
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: shift the mineral-N column headers left by one (drop the old
# "ProfileN" label) and rename them from the old SoilNitrogen naming to the
# new Nutrient naming, adding a 7th column.
$ws.Range("C1").Value = "Soil.Nutrient.MineralN(1)"
$ws.Range("D1").Value = "Soil.Nutrient.MineralN(2)"
$ws.Range("E1").Value = "Soil.Nutrient.MineralN(3)"
$ws.Range("F1").Value = "Soil.Nutrient.MineralN(4)"
$ws.Range("G1").Value = "Soil.Nutrient.MineralN(5)"
$ws.Range("H1").Value = "Soil.Nutrient.MineralN(6)"
$ws.Range("I1").Value = "Soil.Nutrient.MineralN(7)"

# Update the active selection to match (was H1, now I1).
$ws.Range("I1").Select()
